$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1221161.2
$ws.Range("I39").Value = 1587419.6
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 4762258.800000001
$ws.Range("L39").Value = 900
$ws.Range("M39").Value = -4761962.800000001
$ws.Range("N39").Value = -1492

$ws.Range("H69").Value = 8526.666999999999
$ws.Range("I69").Value = 7000
$ws.Range("J69").Value = 9290
$ws.Range("K69").Value = 21000
$ws.Range("L69").Value = 27870
$ws.Range("M69").Value = -20126
$ws.Range("N69").Value = -29618

$ws.Range("H72").Value = 8526.666999999999
$ws.Range("I72").Value = 7000
$ws.Range("J72").Value = 9290
$ws.Range("K72").Value = 63000
$ws.Range("L72").Value = 83610
$ws.Range("M72").Value = -58632
$ws.Range("N72").Value = -92346

$ws.Range("H76").Value = 3875
$ws.Range("I76").Value = 3200
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 3200
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -2885
$ws.Range("N76").Value = -5630

$ws.Range("H79").Value = 3875
$ws.Range("I79").Value = 3200
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 3200
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -2108
$ws.Range("N79").Value = -7184

$ws.Range("H86").Value = 3767
$ws.Range("I86").Value = 2312.6875
$ws.Range("K86").Value = 2312.6875
$ws.Range("M86").Value = -1189.6875

$ws.Range("H89").Value = 3767
$ws.Range("I89").Value = 2312.6875
$ws.Range("K89").Value = 11563.4375
$ws.Range("M89").Value = -5947.4375

$ws.Range("H100").Value = 2061.1667
$ws.Range("I100").Value = 1985.5
$ws.Range("J100").Value = 2212.5
$ws.Range("K100").Value = 1985.5
$ws.Range("L100").Value = 2212.5
$ws.Range("M100").Value = -1444.5
$ws.Range("N100").Value = -3294.5

$ws.Range("H112").Value = 1086.75
$ws.Range("J112").Value = 1106.2858
$ws.Range("L112").Value = 3318.8574
$ws.Range("N112").Value = -5534.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21607.828
$ws.Range("I32").Value = 3476.2576
$ws.Range("J32").Value = 101386.734
$ws.Range("K32").Value = 3476.2576
$ws.Range("L32").Value = 101386.734
$ws.Range("M32").Value = -3189.2576
$ws.Range("N32").Value = -101960.734

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 250749.75
$ws.Range("I94").Value = 250749.75
$ws.Range("K94").Value = 250749.75
$ws.Range("M94").Value = -250298.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16706.121
$ws.Range("I31").Value = 28748.223
$ws.Range("J31").Value = 2255.6
$ws.Range("K31").Value = 28748.223
$ws.Range("L31").Value = 2255.6
$ws.Range("M31").Value = -28453.223
$ws.Range("N31").Value = -2845.6

$ws.Range("H34").Value = 16706.121
$ws.Range("I34").Value = 28748.223
$ws.Range("J34").Value = 2255.6
$ws.Range("K34").Value = 28748.223
$ws.Range("L34").Value = 2255.6
$ws.Range("M34").Value = -28546.223
$ws.Range("N34").Value = -2659.6

$ws.Range("H58").Value = 10117.267
$ws.Range("I58").Value = 1509.9445
$ws.Range("J58").Value = 23028.25
$ws.Range("K58").Value = 1509.9445
$ws.Range("L58").Value = 23028.25
$ws.Range("M58").Value = -1306.9445
$ws.Range("N58").Value = -23434.25

$ws.Range("H134").Value = 1872.6923
$ws.Range("I134").Value = 1302.909
$ws.Range("K134").Value = 3908.727
$ws.Range("M134").Value = -1373.727

$ws.Range("H136").Value = 10117.267
$ws.Range("I136").Value = 1509.9445
$ws.Range("J136").Value = 23028.25
$ws.Range("K136").Value = 4529.833500000001
$ws.Range("L136").Value = 69084.75
$ws.Range("M136").Value = -1979.833500000001
$ws.Range("N136").Value = -74184.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6639
$ws.Range("I5").Value = 580.4783
$ws.Range("J5").Value = 19306.818
$ws.Range("K5").Value = 1741.4349
$ws.Range("L5").Value = 57920.454
$ws.Range("M5").Value = -1629.4349
$ws.Range("N5").Value = -58144.454

$ws.Range("H34").Value = 1833.3334
$ws.Range("J34").Value = 2675
$ws.Range("L34").Value = 8025
$ws.Range("N34").Value = -8193

$ws.Range("H68").Value = 1788.6154
$ws.Range("I68").Value = 921.5909
$ws.Range("J68").Value = 2232.2092
$ws.Range("K68").Value = 2764.7727
$ws.Range("L68").Value = 6696.6276
$ws.Range("M68").Value = -1953.7727
$ws.Range("N68").Value = -8318.6276

$ws.Range("H71").Value = 1788.6154
$ws.Range("I71").Value = 921.5909
$ws.Range("J71").Value = 2232.2092
$ws.Range("K71").Value = 8294.3181
$ws.Range("L71").Value = 20089.8828
$ws.Range("M71").Value = -4238.3181
$ws.Range("N71").Value = -28201.8828

$ws.Range("H122").Value = 8520.154
$ws.Range("I122").Value = 254.66667
$ws.Range("K122").Value = 2292.00003
$ws.Range("M122").Value = 157.9999699999998

$ws.Range("H135").Value = 6639
$ws.Range("I135").Value = 580.4783
$ws.Range("J135").Value = 19306.818
$ws.Range("K135").Value = 5224.3047
$ws.Range("L135").Value = 173761.362
$ws.Range("M135").Value = -2689.3047
$ws.Range("N135").Value = -178831.362

$ws.Range("H139").Value = 1821.1666
$ws.Range("I139").Value = 827.875
$ws.Range("J139").Value = 3807.75
$ws.Range("K139").Value = 2483.625
$ws.Range("L139").Value = 11423.25
$ws.Range("M139").Value = 2656.375
$ws.Range("N139").Value = -21703.25

$ws.Range("H141").Value = 2488.3845
$ws.Range("I141").Value = 1959
$ws.Range("J141").Value = 5400
$ws.Range("K141").Value = 5877
$ws.Range("L141").Value = 16200
$ws.Range("M141").Value = -697
$ws.Range("N141").Value = -26560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 15725
$ws.Range("J6").Value = 15725
$ws.Range("L6").Value = 15725
$ws.Range("N6").Value = -15951

$ws.Range("H16").Value = 15725
$ws.Range("J16").Value = 15725
$ws.Range("L16").Value = 15725
$ws.Range("N16").Value = -16225

$ws.Range("H43").Value = 4189
$ws.Range("I43").Value = 1950
$ws.Range("J43").Value = 5184.1113
$ws.Range("K43").Value = 1950
$ws.Range("L43").Value = 5184.1113
$ws.Range("M43").Value = -1799
$ws.Range("N43").Value = -5486.1113

$ws.Range("H46").Value = 10266.667
$ws.Range("J46").Value = 10266.667
$ws.Range("L46").Value = 10266.667
$ws.Range("N46").Value = -10578.667

$ws.Range("H57").Value = 11520
$ws.Range("J57").Value = 18866.666
$ws.Range("L57").Value = 18866.666
$ws.Range("N57").Value = -20506.666

$ws.Range("H58").Value = 13390
$ws.Range("I58").Value = 1975
$ws.Range("J58").Value = 21000
$ws.Range("K58").Value = 1975
$ws.Range("L58").Value = 21000
$ws.Range("M58").Value = -1698
$ws.Range("N58").Value = -21554

$ws.Range("H70").Value = 41082.465
$ws.Range("I70").Value = 67770.875
$ws.Range("J70").Value = 5497.9165
$ws.Range("K70").Value = 67770.875
$ws.Range("L70").Value = 5497.9165
$ws.Range("M70").Value = -67500.875
$ws.Range("N70").Value = -6037.9165

$ws.Range("H73").Value = 41082.465
$ws.Range("I73").Value = 67770.875
$ws.Range("J73").Value = 5497.9165
$ws.Range("K73").Value = 67770.875
$ws.Range("L73").Value = 5497.9165
$ws.Range("M73").Value = -66834.875
$ws.Range("N73").Value = -7369.9165

$ws.Range("H80").Value = 3825
$ws.Range("J80").Value = 2166.6667
$ws.Range("L80").Value = 2166.6667
$ws.Range("N80").Value = -4162.6667

$ws.Range("H83").Value = 3825
$ws.Range("J83").Value = 2166.6667
$ws.Range("L83").Value = 10833.3335
$ws.Range("N83").Value = -20817.3335

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 2627.8147
$ws.Range("I132").Value = 2734.6667
$ws.Range("K132").Value = 8204.000100000001
$ws.Range("M132").Value = -5674.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1828.9333
$ws.Range("I7").Value = 1362.8334
$ws.Range("J7").Value = 3693.3333
$ws.Range("K7").Value = 1362.8334
$ws.Range("L7").Value = 3693.3333
$ws.Range("M7").Value = -1250.8334
$ws.Range("N7").Value = -3917.3333

$ws.Range("H43").Value = 6837.769
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 6837.769
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 6837.769
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -7223.769

$ws.Range("H93").Value = 1315.2354
$ws.Range("I93").Value = 1304.3077
$ws.Range("J93").Value = 1350.75
$ws.Range("K93").Value = 1304.3077
$ws.Range("L93").Value = 1350.75
$ws.Range("M93").Value = -56.30770000000007
$ws.Range("N93").Value = -3846.75

$ws.Range("H126").Value = 1828.9333
$ws.Range("I126").Value = 1362.8334
$ws.Range("J126").Value = 3693.3333
$ws.Range("K126").Value = 4088.5002
$ws.Range("L126").Value = 11079.9999
$ws.Range("M126").Value = -1618.5002
$ws.Range("N126").Value = -16019.9999

$ws.Range("H132").Value = 6180.421
$ws.Range("I132").Value = 6180.421
$ws.Range("K132").Value = 18541.263
$ws.Range("M132").Value = -16011.263

$ws.Range("H136").Value = 1137.1515
$ws.Range("I136").Value = 1009.46155
$ws.Range("J136").Value = 1611.4286
$ws.Range("K136").Value = 3028.38465
$ws.Range("L136").Value = 4834.2858
$ws.Range("M136").Value = -478.38465
$ws.Range("N136").Value = -9934.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6695.8335
$ws.Range("J74").Value = 7475.6665
$ws.Range("L74").Value = 7475.6665
$ws.Range("N74").Value = -9347.666499999999

$ws.Range("H77").Value = 6695.8335
$ws.Range("J77").Value = 7475.6665
$ws.Range("L77").Value = 22426.9995
$ws.Range("N77").Value = -31786.9995

$ws.Range("H126").Value = 1628.6471
$ws.Range("I126").Value = 1730.1666
$ws.Range("J126").Value = 1385
$ws.Range("K126").Value = 5190.4998
$ws.Range("L126").Value = 4155
$ws.Range("M126").Value = -2720.4998
$ws.Range("N126").Value = -9095

$ws.Range("H132").Value = 3730.1562
$ws.Range("I132").Value = 4973.0586
$ws.Range("J132").Value = 2321.5334
$ws.Range("K132").Value = 14919.1758
$ws.Range("L132").Value = 6964.600199999999
$ws.Range("M132").Value = -12389.1758
$ws.Range("N132").Value = -12024.6002

$ws.Range("H136").Value = 715.2292
$ws.Range("I136").Value = 473.08823
$ws.Range("J136").Value = 1303.2858
$ws.Range("K136").Value = 1419.26469
$ws.Range("L136").Value = 3909.8574
$ws.Range("M136").Value = 1130.73531
$ws.Range("N136").Value = -9009.857400000001
